# "Generate Report for Archive"
# 1) All cells whose status was "Ready for handoff" become "In Translation"
#    (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all share the same string).
# 2) The "Status" columns get narrower: Overview columns E & F, and the
#    Status column (C) on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# --- 1. Update status text everywhere it appears ---------------------------
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# --- 2. Narrow the status columns -------------------------------------------
# Target stored width is 13.4101848602295 character-units; the COM layer
# snaps ColumnWidth to whole pixels, so feed it the de-pixelated value that
# lands closest to that target.
$newColumnWidth = 12.576851526896165

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
